# "start of fork controls"
#
# The hardware map has several rows reserved (PORT 2/A-H) whose DEVICE
# cell (column B) was still an empty, strikethrough "placeholder" cell.
# This change starts wiring up the new fork subsystem:
#   - PORT 2  -> inertial sensor (IMU)
#   - PORT E  -> fork limit switch
#   - PORT F/G -> the two mobile-goal-lock pneumatics solenoids
#
# Cell writes are ordered to match the author's original editing session
# (NAME/FUNCTION text typed in before the DEVICE column for the new
# rows), so the shared-string table comes out in the same order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# PORT 2 (row 3) - NAME / FUNCTION first
$ws.Range("C3").Value = "imu"
$ws.Range("D3").Value = "inertial sensor"

# PORT E (row 26) - NAME / FUNCTION first
$ws.Range("C26").Value = "lim"
$ws.Range("D26").Value = "fork limit switch"

# PORT F (row 27) - DEVICE first
$ws.Range("B27").Value = "Pneumatics*"

# PORT E (row 26) - DEVICE, filled in after F's DEVICE
$ws.Range("B26").Value = "Limit*"

# PORT F / G (rows 27, 28) - NAME
$ws.Range("C27").Value = "mogo_lock1"
$ws.Range("C28").Value = "mogo_lock2"

# PORT 2 (row 3) - DEVICE, filled in last
$ws.Range("B3").Value = "Inertial*"

# PORT G (row 28) - DEVICE (reuses the "Pneumatics*" shared string)
$ws.Range("B28").Value = "Pneumatics*"

# These two cells were blank "reserved" placeholders styled with strike-
# through text; now that they hold real data, clear the strikethrough so
# they read as normal, filled-in rows.
$ws.Range("B3").Font.Strikethrough = $false
$ws.Range("B26").Font.Strikethrough = $false

# Match the author's final viewport / selection state.
$ws.Range("B3").Select()
